# Initial importation of VL / SL (vacation leave / sick leave) columns into
# the "Pool" sheet, appended right after the existing a_1..a_12 activity
# columns (AO). Not yet referenced anywhere else in the workbook (matches
# the commit message: "Not yet used in object assignment").

$wb = $excel.ActiveWorkbook
$pool = $wb.Worksheets.Item("Pool")

# Header row (row 1): AP1 = "VL", AQ1 = "SL"
$pool.Range("AP1").Value = "VL"
$pool.Range("AQ1").Value = "SL"

# Data row (row 2): AP2 = AQ2 = "15" (stored as text, like the neighbouring
# text-formatted cells in that row/column block)
$pool.Range("AP2").Value = "15"
$pool.Range("AQ2").Value = "15"

# Leave the new range selected/active on the Pool sheet, like Excel would
# right after typing the new values in (also keeps Pool as the active tab).
$pool.Range("AP1:AQ2").Select() | Out-Null
